# Week 17 data log + numeric roll-up updates across OFF/DEF/ST/TURNS/PEN sheets
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 17 per-game yardage logs (space-separated run) ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value() + " 6 0 1 1 3 4 3 20 1 16 6 6 4 12 3 5 20 9 5 -2 14 4 7 3 6 11 5 5 4 2 2 12 4 4 6 5 4 2 3 8"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 13 15 11 15 5 6 10 13 9 15 8 2"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 2 4 3 0 4 4 13 9 -4 8 6 4 1 11 -3 7 5 2 13 0 6"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 8 15 0 9 10 61 15 24 4 1 8 15 17 11"

# --- OFF sheet: updated season totals after Week 17 ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 333
$ws.Range("D2").Value = 24
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 145
$ws.Range("G2").Value = 132
$ws.Range("J2").Value = 74
$ws.Range("L2").Value = 631
$ws.Range("M2").Value = 406
$ws.Range("O2").Value = 33
$ws.Range("P2").Value = 21
$ws.Range("Q2").Value = 1115
$ws.Range("C3").Value = 535
$ws.Range("D3").Value = 19
$ws.Range("E3").Value = 69
$ws.Range("F3").Value = 224
$ws.Range("G3").Value = 80
$ws.Range("I3").Value = 107
$ws.Range("J3").Value = 118
$ws.Range("N3").Value = 36

# --- DEF sheet: updated season totals after Week 17 ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 427
$ws.Range("D2").Value = 20
$ws.Range("F2").Value = 108
$ws.Range("G2").Value = 135
$ws.Range("L2").Value = 588
$ws.Range("M2").Value = 350
$ws.Range("O2").Value = 58
$ws.Range("Q2").Value = 1114
$ws.Range("C3").Value = 390
$ws.Range("E3").Value = 85
$ws.Range("F3").Value = 217
$ws.Range("G3").Value = 87
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 132
$ws.Range("J3").Value = 99
$ws.Range("N3").Value = 39

# --- ST sheet: updated season totals + appended Week 17 per-game logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 205
$ws.Range("F2").Value = 107
$ws.Range("G2").Value = 105
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 3
$ws.Range("B3").Value = 130
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 64"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 23"
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 23 26 25"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 0 0 0 0 19"

# --- TURNS sheet: updated season totals after Week 17 ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 13
$ws.Range("E2").Value = 26
$ws.Range("E3").Value = 25

# --- PEN sheet: updated season totals after Week 17 ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 21
$ws.Range("D4").Value = 17
